$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing week 22 value
$ws.Range("B23").Value = 346

# Add week 23, 24, 25 data
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 357

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 261

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 51
